$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "29.790.61"
$ws.Range("E2").Value2 = "  -0.43%  "
$ws.Range("D3").Value2 = "1.869.91"
$ws.Range("E3").Value2 = "  -0.27%  "
$ws.Range("E4").Value2 = "  -0.02%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value2 = "0.7267"
$c.Style = "Normal"
$ws.Range("E5").Value2 = "  -1.99%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value2 = "241.14"
$c.Style = "Normal"
$ws.Range("E6").Value2 = "  -0.54%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value2 = "1.001"
$c.Style = "Normal"
$ws.Range("E7").Value2 = "  +0.03%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value2 = "0.3130"
$c.Style = "Normal"
$ws.Range("E8").Value2 = "  -0.75%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value2 = "0.07133"
$c.Style = "Normal"
$ws.Range("E9").Value2 = "  -0.62%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value2 = "24.38"
$c.Style = "Normal"
$ws.Range("E10").Value2 = "  -1.52%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value2 = "0.08124"
$c.Style = "Normal"
$ws.Range("E11").Value2 = "  -3.90%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value2 = "0.7419"
$c.Style = "Normal"
$ws.Range("E12").Value2 = "  -1.56%  "
$ws.Range("D13").Value2 = "1.882.47"
$ws.Range("E13").Value2 = "  -0.04%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value2 = "5.341"
$c.Style = "Normal"
$ws.Range("E14").Value2 = "  -1.11%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value2 = "92.42"
$c.Style = "Normal"
$ws.Range("E15").Value2 = "  -0.19%  "
$ws.Range("D16").Value2 = "29.793.26"
$ws.Range("E16").Value2 = "  -0.45%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value2 = "5.995"
$c.Style = "Normal"
$ws.Range("E17").Value2 = "  -1.81%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value2 = "247.19"
$c.Style = "Normal"
$ws.Range("E18").Value2 = "  +1.50%  "
$ws.Range("E19").Value2 = "  -1.68%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value2 = "0.000007798"
$c.Style = "Normal"
$ws.Range("E20").Value2 = "  -0.32%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value2 = "0.9999"
$c.Style = "Normal"
$ws.Range("E21").Value2 = "  +0.01%  "
$ws.Range("D22").Value2 = "2.127.69"
$ws.Range("E22").Value2 = "  +0.22%  "
$ws.Range("E23").Value2 = "  +0.03%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value2 = "7.736"
$c.Style = "Normal"
$ws.Range("E24").Value2 = "  -3.31%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value2 = "0.1531"
$c.Style = "Normal"
$ws.Range("E25").Value2 = "  -1.89%  "
$ws.Range("E26").Value2 = "  -1.35%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value2 = "163.66"
$c.Style = "Normal"
$ws.Range("E27").Value2 = "  -1.29%  "
$ws.Range("E28").Value2 = "  -0.59%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value2 = "2.001"
$c.Style = "Normal"
$ws.Range("E29").Value2 = "  -2.11%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value2 = "1.445"
$c.Style = "Normal"
$ws.Range("E30").Value2 = "  -2.27%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value2 = "4.517"
$c.Style = "Normal"
$ws.Range("E31").Value2 = "  -2.03%  "
$ws.Range("E32").Value2 = "  -0.33%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value2 = "4.172"
$c.Style = "Normal"
$ws.Range("E33").Value2 = "  -2.51%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value2 = "0.05316"
$c.Style = "Normal"
$ws.Range("E34").Value2 = "  -0.47%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value2 = "1.226"
$c.Style = "Normal"
$ws.Range("E35").Value2 = "  -1.46%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value2 = "0.7370"
$c.Style = "Normal"
$ws.Range("E36").Value2 = "  -2.67%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value2 = "0.9979"
$c.Style = "Normal"
$ws.Range("E37").Value2 = "  -0.21%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value2 = "2.697"
$c.Style = "Normal"
$ws.Range("E38").Value2 = "  -0.02%  "
$ws.Range("E39").Value2 = "  -1.25%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value2 = "2.734"
$c.Style = "Normal"
$ws.Range("E40").Value2 = "  -0.64%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value2 = "0.4468"
$c.Style = "Normal"
$ws.Range("E41").Value2 = "  -0.48%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value2 = "0.8853"
$c.Style = "Normal"
$ws.Range("E42").Value2 = "  +2.97%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value2 = "5.962"
$c.Style = "Normal"
$ws.Range("E43").Value2 = "  -2.30%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value2 = "71.18"
$c.Style = "Normal"
$ws.Range("D45").Value2 = "1.041.23"
$ws.Range("E45").Value2 = "  -6.51%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value2 = "1.001"
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value2 = "103.77"
$c.Style = "Normal"
$ws.Range("E47").Value2 = "  +0.56%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value2 = "1.820"
$c.Style = "Normal"
$ws.Range("E48").Value2 = "  -1.19%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value2 = "7.454"
$c.Style = "Normal"
$ws.Range("E49").Value2 = "  -3.03%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value2 = "9.535"
$c.Style = "Normal"
$ws.Range("E50").Value2 = "  -0.29%  "
$ws.Range("D51").Value2 = "2.019.38"
$ws.Range("E51").Value2 = "  -0.16%  "
